$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line entries ("line7", "line8") are being inserted into the lines
# table right before the "extr1.." block (i.e. at row 8), pushing the
# existing extr1..extr8 rows down by two rows (old rows 8-15 -> new rows
# 10-17). Shift the name/from_bus/to_bus/in_service columns (B:E) for those
# rows down first (bottom-up so we don't clobber data we still need), then
# fill in the two freshly-opened rows with the new line data.
for ($r = 15; $r -ge 8; $r--) {
    $dest = $r + 2
    $ws.Range("B$r`:E$r").Copy($ws.Range("B$dest`:E$dest"))
}

# Column A is just the 0-based row counter - recompute it for the rows that
# moved down. Rows 16/17 are brand new cells, so first clone the bold/
# bordered "index column" formatting from an existing A cell before writing
# the value into them.
$ws.Range("A2").Copy($ws.Range("A16"))
$ws.Range("A2").Copy($ws.Range("A17"))
for ($r = 10; $r -le 17; $r++) {
    $ws.Range("A$r").Value = $r - 2
}

# New row for line7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $false

# New row for line8
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# extr1 and extr2 (now at rows 10 and 11) flip to in_service = TRUE
$ws.Range("E10").Value = $true
$ws.Range("E11").Value = $true
